# "Add files via upload" — bets-2023-2.xlsx
#
# Net effect of the authored edit (per the OOXML diff):
#   - Column C (the "ex_date" helper column) on the "bets" sheet held a
#     shared formula =TEXT(Bn,"YYYY-MM-DD"). It was converted to plain,
#     static text values (same visible dates), which is what you get from
#     a copy / paste-special-values over that range in Excel. This also
#     ripples the remaining shared-formula si= indices down by one and
#     grows the shared-strings table with the 23 distinct date strings.
#   - The sheet's view scrolled back to the top and the selection moved
#     to E4 (no more frozen/scrolled "A13" top-left, no more F44 selection).
#   - Column C picked up an explicit (best-fit-ish) column width.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

# --- 1. Freeze column C's TEXT() formulas into static values -------------
# Copy / Paste-Special-Values in place: keeps the existing cell style
# (so the date still displays as m/d/yyyy-ish text) while turning every
# formula cell into a literal value, exactly like the authored diff shows
# (t="str" + <f> -> t="s" shared-string literal).
$colC = $ws.Range("C2:C43")
$colC.Copy() | Out-Null
$colC.PasteSpecial(-4163) | Out-Null   # xlPasteValues
$excel.CutCopyMode = 0

# --- 2. Reset the view: scroll to top, select E4 --------------------------
$ws.Range("E4").Select() | Out-Null

# --- 3. Give column C its own width now that it holds fitted text --------
$ws.Columns("C").ColumnWidth = 9.2
